$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the "type" column (C) values that were not correctly taken care of
# during a previous edit.
$ws.Range("C1").Value = "EC"
$ws.Range("C4").Value = "IUT"
$ws.Range("C5").Value = "EI"

# F4 holds the numeric-looking establishment code "7764", but (like F1/F8)
# it should be stored as text rather than a number. Copy the already-text
# value from F1 (same underlying text "7764") using Paste Values so the
# cell keeps its own existing style/format instead of picking up a new one.
$ws.Range("F1").Copy()
$ws.Range("F4").PasteSpecial(-4163)
